# Hortaliza, Macroferia Regional de Talca - Sandia
# Insert a new weekly price record at row 169 (above the existing 2021-01-07
# "Extra" record), pushing the existing rows 169:190 down to 170:191.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(169).Insert()

$ws.Cells.Item(169, 1).Value = 5
$ws.Cells.Item(169, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(169, 3).Value = 'Maule'
$ws.Cells.Item(169, 4).Value = 44900
$ws.Cells.Item(169, 5).Value = 7
$ws.Cells.Item(169, 6).Value = 100112028
$ws.Cells.Item(169, 7).Value = 'Sandia'
$ws.Cells.Item(169, 8).Value = 'Sin especificar'
$ws.Cells.Item(169, 9).Value = 'Primera'
$ws.Cells.Item(169, 10).Value = 3000
$ws.Cells.Item(169, 11).Value = 400
$ws.Cells.Item(169, 12).Value = 400
$ws.Cells.Item(169, 13).Value = 400
$ws.Cells.Item(169, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(169, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(169, 16).Value = 400
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = 'Hortaliza'
